# Updates cryptos list Price (D) and Volume(1h) (E) columns for Sheet1.
# Values that look purely numeric are prefixed with a leading apostrophe so
# Excel keeps them as text (matching the original inline-string cells)
# instead of silently re-parsing/rounding them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.714.52'
$ws.Range("E2").Value = '  +2.17%  '

$ws.Range("D3").Value = '2.941.37'
$ws.Range("E3").Value = '  +0.41%  '

$ws.Range("D5").Value = '''592.87'
$ws.Range("E5").Value = '  -0.83%  '

$ws.Range("D6").Value = '''147.25'
$ws.Range("E6").Value = '  +1.34%  '

$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").Value = '2.939.08'
$ws.Range("E8").Value = '  +0.35%  '

$ws.Range("E9").Value = '  +0.64%  '

$ws.Range("D10").Value = '''7.30'
$ws.Range("E10").Value = '  +4.73%  '

$ws.Range("D11").Value = '''0.149'
$ws.Range("E11").Value = '  +5.13%  '

$ws.Range("D12").Value = '''0.439'
$ws.Range("E12").Value = '  +0.24%  '

$ws.Range("D13").Value = '''0.0000235'
$ws.Range("E13").Value = '  +4.53%  '

$ws.Range("D14").Value = '''32.56'
$ws.Range("E14").Value = '  -2.79%  '

$ws.Range("E15").Value = '  -0.79%  '

$ws.Range("D16").Value = '3.429.12'
$ws.Range("E16").Value = '  +0.44%  '

$ws.Range("D17").Value = '62.686.28'
$ws.Range("E17").Value = '  +2.10%  '

$ws.Range("D18").Value = '''6.68'
$ws.Range("E18").Value = '  -0.24%  '

$ws.Range("D19").Value = '2.944.38'
$ws.Range("E19").Value = '  +0.45%  '

$ws.Range("D20").Value = '''439.22'
$ws.Range("E20").Value = '  +1.81%  '

$ws.Range("E21").Value = '  -0.96%  '

$ws.Range("D22").Value = '''0.663'
$ws.Range("E22").Value = '  -1.92%  '

$ws.Range("D23").Value = '''7.01'
$ws.Range("E23").Value = '  -0.82%  '

$ws.Range("D24").Value = '''80.82'
$ws.Range("E24").Value = '  -1.25%  '

$ws.Range("D25").Value = '''11.14'
$ws.Range("E25").Value = '  +2.55%  '

$ws.Range("E26").Value = '  -2.70%  '

$ws.Range("E27").Value = '  -0.52%  '

$ws.Range("E28").Value = '  -0.05%  '

$ws.Range("E29").Value = '  +0.66%  '

$ws.Range("E30").Value = '  +3.47%  '

$ws.Range("E31").Value = '  -0.38%  '

$ws.Range("E32").Value = '  +14.36%  '

$ws.Range("E33").Value = '  -1.32%  '

$ws.Range("E34").Value = '  -1.20%  '

$ws.Range("E35").Value = '  +0.00%  '

$ws.Range("D36").Value = '''0.991'
$ws.Range("E36").Value = '  -2.20%  '

$ws.Range("E37").Value = '  +2.97%  '

$ws.Range("D38").Value = '''5.56'
$ws.Range("E38").Value = '  -1.12%  '

$ws.Range("D39").Value = '''49.66'
$ws.Range("E39").Value = '  -0.62%  '

$ws.Range("E40").Value = '  +1.00%  '

$ws.Range("E41").Value = '  -1.40%  '

$ws.Range("E42").Value = '  -4.30%  '

$ws.Range("E43").Value = '  -1.01%  '

$ws.Range("D44").Value = '''38.65'
$ws.Range("E44").Value = '  -8.86%  '

$ws.Range("D45").Value = '2.692.38'
$ws.Range("E45").Value = '  -0.21%  '

$ws.Range("D46").Value = '''134.65'
$ws.Range("E46").Value = '  +0.78%  '

$ws.Range("D47").Value = '''359.78'
$ws.Range("E47").Value = '  -0.82%  '

$ws.Range("E48").Value = '  -3.46%  '

$ws.Range("E50").Value = '  -1.10%  '

$ws.Range("E51").Value = '  -4.36%  '
